$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. The cached "today" date field ("Espaço Reservado para Data") shows
#    19/02/2026 on the Slide Master and on every Custom Layout; bump it
#    to 26/02/2026 everywhere it appears.
# ---------------------------------------------------------------------
$oldDate = "19/02/2026"
$newDate = "26/02/2026"

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

$layouts = $master.CustomLayouts
for ($l = 1; $l -le $layouts.Count; $l++) {
    $layout = $layouts.Item($l)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2. Slide 17 ("Exercícios") — the box "CaixaDeTexto 10":
#      * merge the first bullet's 3 runs into one clean sentence (the
#        stray "Por que não é recomendável usar..." duplicate question
#        is removed),
#      * append a brand-new bullet about FK/dependents,
#      * append a trailing empty bullet paragraph,
#      * grow the shape so the extra bullet still fits.
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(17)

$target = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $candidate = $slide.Shapes.Item($i)
    if ($candidate.Name -eq "CaixaDeTexto 10") {
        $target = $candidate
    }
}

$bullet1 = 'Se você precisa guardar apenas "Sim" ou "Não", por que usar um campo de texto livre (VARCHAR) seria uma má escolha de modelagem? '
$bullet2 = 'No meio de uma transferência bancária, a energia do servidor cai. Por que o saldo não pode ter saído de uma conta sem ter entrado na outra?'
$bullet3 = 'Por que as Chaves Estrangeiras são vitais para evitar que o banco de dados se torne um amontoado de "dados órfãos" (informações que existem, mas não se ligam a nada)?'
$bullet4 = 'Pense no contexto de uma empresa, um empregado tem um dependente (filho), qual a relação da empresa com este dependente? O que aconteceria com o dependente se o empregado fosse demitido? Dentro da relação entre empregado e dependente, em qual das entidades ficaria a chave estrangeira? Por que?'

$target.TextFrame.TextRange.Text = $bullet1 + "`r" + $bullet2 + "`r" + $bullet3 + "`r" + $bullet4 + "`r"

# Resize the textbox (width/offset stay put, only the height grows) —
# this must happen *after* the text is set, otherwise the shape's
# auto-fit recomputes Height from the new text and clobbers our value.
$target.Height = 298.0828552246094
